# edit.ps1 - applies the Admin.docx change:
#  1. Strip the <w:proofErr .../> spell-check markers that bracket every
#     "...@...local" style login/email paragraph.
#  2. Remove the _GoBack bookmark from the trailing "php.exe" command
#     paragraph, append three new paragraphs (blank / Arabic "بقوقل" /
#     "localhost:8000") after it, move the _GoBack bookmark onto the new
#     last paragraph, and drop the now-superfluous empty paragraph that
#     used to sit right before the sectPr.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove proofErr spellStart/spellEnd wrappers around the email
# / login paragraphs. Round-tripping a paragraph's Range through
# XML()/InsertXML() rebuilds the paragraph from the canonical OOXML
# model, which does not retain proofErr markers, while it keeps every
# other bit of paragraph/run formatting untouched.
# ---------------------------------------------------------------------
$targets = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "@.*\.local") {
        $targets += $p
    }
}
foreach ($p in $targets) {
    $r = $p.Range
    $xml = $r.XML()
    $r.InsertXML($xml)
}

# ---------------------------------------------------------------------
# Step 2: rebuild the document tail starting at the paragraph that holds
# the php.exe command line through the end of the document (this also
# covers the trailing empty paragraph and the sectPr). Doing this as one
# InsertXML keeps the bookmark relocation, the new paragraphs and the
# deletion of the trailing empty paragraph atomic and well formed.
# ---------------------------------------------------------------------
$phpPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*php.exe*") {
        $phpPara = $p
    }
}

$tailStart = $phpPara.Range.Start
$tailEnd = $d.Content.End
$tailRange = $d.Range($tailStart, $tailEnd)

$php = '&amp; &quot;C:\xampp\php\php.exe&quot; -S localhost:8000 -t &quot;C:\xampp\htdocs\Traveler-Companion&quot;'
$arabic = "بقوقل"

$body = ""
$body += "<w:p><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:t>" + $php + "</w:t></w:r></w:p>"
$body += "<w:p><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr></w:p>"
$body += "<w:p><w:pPr><w:rPr><w:rtl/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint=`"cs`"/><w:rtl/></w:rPr><w:t>" + $arabic + "</w:t></w:r></w:p>"
$body += "<w:p><w:pPr><w:rPr><w:rFonts w:hint=`"cs`"/></w:rPr></w:pPr><w:r><w:t>localhost:8000</w:t></w:r><w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/></w:p>"
$body += "<w:sectPr w:rsidR=`"00317B67`" w:rsidRPr=`"00317B67`" w:rsidSect=`"00986D99`"><w:pgSz w:w=`"11906`" w:h=`"16838`"/><w:pgMar w:top=`"1440`" w:right=`"1800`" w:bottom=`"1440`" w:left=`"1800`" w:header=`"708`" w:footer=`"708`" w:gutter=`"0`"/><w:cols w:space=`"708`"/><w:bidi/><w:rtlGutter/><w:docGrid w:linePitch=`"360`"/></w:sectPr>"

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $body + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$tailRange.InsertXML($xml)
